$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B5").Value = "Ratio with UCUM or EDQM codes if code is used"
$wsMeta.Range("B8").Value = "2025-08-13T14:10:49+00:00"
$wsMeta.Range("B12").Value = "Ratio with numerator and denominator unit UCUM or EDQM encoded if code is used"

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("L5").Value = "Numerator value"
$wsElem.Range("M5").Value = "The value of the numerator."
$wsElem.Range("L6").Value = "Denominator value"
$wsElem.Range("M6").Value = "The value of the denominator."
